$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Fecha (date) + price columns (Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# for rows 2-21. The underlying date/price rows are reshuffled among the existing records.
$updates = @{
    2 = @{ D = 44638; J = 1000; K = 900; L = 950; M = 925; P = 925 }
    3 = @{ D = 44455; J = 1100; K = 900; L = 1000; M = 950; P = 950 }
    4 = @{ D = 44284; J = 1500; K = 800; L = 850; M = 825; P = 825 }
    5 = @{ D = 44291; J = 1000; K = 1000; L = 1200; M = 1100; P = 1100 }
    6 = @{ D = 44407; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
    7 = @{ D = 44649; J = 600; K = 900; L = 1000; M = 950; P = 950 }
    8 = @{ D = 44341; J = 1300; K = 900; L = 1000; M = 950; P = 950 }
    9 = @{ D = 44449; J = 1300; K = 900; L = 950; M = 925; P = 925 }
    10 = @{ D = 44442; J = 1250; K = 850; L = 900; M = 875; P = 875 }
    11 = @{ D = 44484; J = 900; K = 750; L = 800; M = 775; P = 775 }
    12 = @{ D = 44175; J = 1600; K = 1000; L = 1200; M = 1100; P = 1100 }
    13 = @{ D = 44550; J = 1300; K = 1000; L = 1200; M = 1100; P = 1100 }
    14 = @{ D = 44607; J = 900; K = 1300; L = 1400; M = 1350; P = 1350 }
    15 = @{ D = 44243; J = 1200; K = 1200; L = 1300; M = 1250; P = 1250 }
    16 = @{ D = 44687; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
    17 = @{ D = 44476; J = 900; K = 700; L = 800; M = 750; P = 750 }
    18 = @{ D = 44453; J = 1000; K = 800; L = 900; M = 850; P = 850 }
    19 = @{ D = 44673; J = 900; K = 1300; L = 1400; M = 1350; P = 1350 }
    20 = @{ D = 44229; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    21 = @{ D = 44656; J = 1000; K = 900; L = 1000; M = 950; P = 950 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
